$wb = $excel.ActiveWorkbook

# ---- Sheet1 (TABLE_1): add 07/01/2022 column (DX), update 06/01/2022 (DW) ----
$ws1 = $wb.Worksheets.Item("TABLE_1")
$ws1.Range("DW4").Copy()
$ws1.Range("DX4").PasteSpecial(-4122)
$ws1.Range("DX4").Value = "'07/01/2022"
$ws1.Range("DW5").Value = 9994.4
$ws1.Range("DX5").Value = 8985.4
$ws1.Range("DW6").Value = 165.8
$ws1.Range("DX6").Value = 157.5
$ws1.Range("DW7").Value = 21.5
$ws1.Range("DX7").Value = 16.6
$ws1.Range("DW8").Value = 150.2
$ws1.Range("DX8").Value = 146.8
$ws1.Range("DW9").Value = 93.9
$ws1.Range("DX9").Value = 82.8
$ws1.Range("DW10").Value = 1223.2
$ws1.Range("DX10").Value = 1053
$ws1.Range("DW11").Value = 195.4
$ws1.Range("DX11").Value = 181.1
$ws1.Range("DW12").Value = 110.4
$ws1.Range("DX12").Value = 97.1
$ws1.Range("DW13").Value = 33.5
$ws1.Range("DX13").Value = 32.7
$ws1.Range("DW14").Copy()
$ws1.Range("DX14").PasteSpecial(-4122)
$ws1.Range("DW15").Value = 399.7
$ws1.Range("DX15").Value = 391
$ws1.Range("DW16").Value = 330.9
$ws1.Range("DX16").Value = 318.9
$ws1.Range("DW17").Copy()
$ws1.Range("DX17").PasteSpecial(-4122)
$ws1.Range("DW18").Value = 60.2
$ws1.Range("DX18").Value = 54.1
$ws1.Range("DW19").Value = 405.1
$ws1.Range("DX19").Value = 377.5
$ws1.Range("DW20").Value = 191.4
$ws1.Range("DX20").Value = 173.1
$ws1.Range("DW21").Value = 131.5
$ws1.Range("DX21").Value = 114
$ws1.Range("DW22").Value = 119.6
$ws1.Range("DX22").Value = 100.9
$ws1.Range("DW23").Value = 147.4
$ws1.Range("DX23").Value = 121.9
$ws1.Range("DW24").Value = 132.8
$ws1.Range("DX24").Value = 127.7
$ws1.Range("DW25").Value = 46.7
$ws1.Range("DX25").Value = 39.7
$ws1.Range("DW26").Value = 202.1
$ws1.Range("DX26").Value = 188.9
$ws1.Range("DW27").Value = 221.5
$ws1.Range("DX27").Value = 194.8
$ws1.Range("DW28").Value = 270.8
$ws1.Range("DX28").Value = 243.4
$ws1.Range("DW29").Value = 192.9
$ws1.Range("DX29").Value = 165.6
$ws1.Range("DW30").Value = 99.2
$ws1.Range("DX30").Value = 93.6
$ws1.Range("DW31").Copy()
$ws1.Range("DX31").PasteSpecial(-4122)
$ws1.Range("DW32").Value = 39.6
$ws1.Range("DX32").Value = 33.3
$ws1.Range("DW33").Value = 88.2
$ws1.Range("DX33").Value = 77.5
$ws1.Range("DW34").Value = 67.7
$ws1.Range("DX34").Value = 64.1
$ws1.Range("DW35").Value = 38.3
$ws1.Range("DX35").Value = 31.7
$ws1.Range("DW36").Value = 302.2
$ws1.Range("DX36").Value = 247.1
$ws1.Range("DW37").Value = 63.8
$ws1.Range("DX37").Value = 61.4
$ws1.Range("DW38").Value = 666.6
$ws1.Range("DX38").Value = 546.9
$ws1.Range("DW39").Value = 288.9
$ws1.Range("DX39").Value = 240
$ws1.Range("DW40").Value = 34.1
$ws1.Range("DX40").Value = 29
$ws1.Range("DW41").Value = 344.9
$ws1.Range("DX41").Value = 330.4
$ws1.Range("DW42").Value = 131.1
$ws1.Range("DX42").Value = 119.9
$ws1.Range("DW43").Value = 129.7
$ws1.Range("DX43").Value = 99.7
$ws1.Range("DW44").Value = 293.7
$ws1.Range("DX44").Value = 264.3
$ws1.Range("DW45").Value = 27.2
$ws1.Range("DX45").Value = 23.6
$ws1.Range("DW46").Value = 159.5
$ws1.Range("DX46").Value = 152.7
$ws1.Range("DW47").Value = 33.2
$ws1.Range("DX47").Value = 28
$ws1.Range("DW48").Value = 191.3
$ws1.Range("DX48").Value = 167.8
$ws1.Range("DW49").Value = 1104.1
$ws1.Range("DX49").Value = 1032.4
$ws1.Range("DW50").Value = 114
$ws1.Range("DX50").Value = 102.2
$ws1.Range("DW51").Value = 28.8
$ws1.Range("DX51").Value = 25.6
$ws1.Range("DW52").Value = 295.2
$ws1.Range("DX52").Value = 261.2
$ws1.Range("DW53").Value = 247.9
$ws1.Range("DX53").Value = 236
$ws1.Range("DW54").Value = 54.3
$ws1.Range("DX54").Value = 49.6
$ws1.Range("DW55").Value = 186.3
$ws1.Range("DX55").Value = 167
$ws1.Range("DW56").Value = 27.3
$ws1.Range("DX56").Value = 24.3

# ---- Sheet2 (TABLE_2): add 07/01/2022 YoY column (DL), update 06/01/2022 YoY (DK) ----
$ws2 = $wb.Worksheets.Item("TABLE_2")
$ws2.Range("DK4").Copy()
$ws2.Range("DL4").PasteSpecial(-4122)
$ws2.Range("DL4").Value = "'07/01/2022"
$ws2.Range("DK5").Value = 1.59181930919515
$ws2.Range("DL5").Value = 1.38902993579545
$ws2.Range("DK6").Value = 2.03076923076924
$ws2.Range("DL6").Value = 1.15606936416186
$ws2.Range("DK7").Value = -8.11965811965811
$ws2.Range("DL7").Value = -17.4129353233831
$ws2.Range("DK8").Value = -0.199335548172765
$ws2.Range("DL8").Value = 0.479123887748129
$ws2.Range("DK9").Value = 2.39912758996729
$ws2.Range("DL9").Value = 1.97044334975369
$ws2.Range("DK10").Value = 6.04247941048982
$ws2.Range("DL10").Value = 4.526503871352
$ws2.Range("DK11").Value = -0.458481915435575
$ws2.Range("DL11").Value = -2.89544235924932
$ws2.Range("DK12").Value = -0.986547085201789
$ws2.Range("DL12").Value = 0.413650465356794
$ws2.Range("DK13").Value = 0.299401197604795
$ws2.Range("DL13").Value = -0.30487804878047
$ws2.Range("DK14").Copy()
$ws2.Range("DL14").PasteSpecial(-4122)
$ws2.Range("DK15").Value = 1.5497967479675
$ws2.Range("DL15").Value = 0.851173587825641
$ws2.Range("DK16").Value = 3.73040752351096
$ws2.Range("DL16").Value = 4.07963446475196
$ws2.Range("DK17").Copy()
$ws2.Range("DL17").PasteSpecial(-4122)
$ws2.Range("DK18").Value = 8.07899461400359
$ws2.Range("DL18").Value = 9.51417004048581
$ws2.Range("DK19").Value = 5.46732621713095
$ws2.Range("DL19").Value = 5.38805136795087
$ws2.Range("DK20").Value = -1.69491525423729
$ws2.Range("DL20").Value = 3.03571428571428
$ws2.Range("DK21").Value = 0.998463901689717
$ws2.Range("DL21").Value = 1.06382978723405
$ws2.Range("DK22").Value = -0.416319733555371
$ws2.Range("DL22").Value = -0.493096646942787
$ws2.Range("DK23").Value = 0.408719346049042
$ws2.Range("DL23").Value = -1.45513338722717
$ws2.Range("DK24").Value = -0.673148840688091
$ws2.Range("DL24").Value = -0.545171339563854
$ws2.Range("DK25").Value = 1.52173913043479
$ws2.Range("DL25").Value = 3.11688311688312
$ws2.Range("DK26").Value = 3.85405960945529
$ws2.Range("DL26").Value = 3.56359649122807
$ws2.Range("DK27").Value = -1.46797153024912
$ws2.Range("DL27").Value = -2.25790265930759
$ws2.Range("DK28").Value = 1.68982350732257
$ws2.Range("DL28").Value = 2.0973154362416
$ws2.Range("DK29").Value = -0.77160493827159
$ws2.Range("DL29").Value = -2.35849056603774
$ws2.Range("DK30").Value = 2.16271884654996
$ws2.Range("DL30").Value = 1.51843817787418
$ws2.Range("DK31").Copy()
$ws2.Range("DL31").PasteSpecial(-4122)
$ws2.Range("DK32").Value = -3.17848410757946
$ws2.Range("DL32").Value = 0.301204819277091
$ws2.Range("DK33").Value = 2.08333333333333
$ws2.Range("DL33").Value = -1.39949109414758
$ws2.Range("DK34").Value = -0.147492625368723
$ws2.Range("DL34").Value = 1.42405063291138
$ws2.Range("DK35").Value = 2.95698924731181
$ws2.Range("DL35").Value = -0.314465408805025
$ws2.Range("DK36").Value = 4.93055555555555
$ws2.Range("DL36").Value = 2.10743801652892
$ws2.Range("DK37").Value = 0.156985871271588
$ws2.Range("DL37").Value = 0.490998363338796
$ws2.Range("DK38").Value = 7.20488903184303
$ws2.Range("DL38").Value = 7.06734534064212
$ws2.Range("DK39").Value = -2.76001346348032
$ws2.Range("DL39").Value = 0.292519849561215
$ws2.Range("DK40").Value = 5.2469135802469
$ws2.Range("DL40").Value = 10.6870229007634
$ws2.Range("DK41").Value = -2.57062146892656
$ws2.Range("DL41").Value = -0.45194335643266
$ws2.Range("DK42").Value = -1.28012048192772
$ws2.Range("DL42").Value = -0.909090909090916
$ws2.Range("DK43").Value = 5.87755102040815
$ws2.Range("DL43").Value = -9.85533453887884
$ws2.Range("DK44").Value = -1.50905432595573
$ws2.Range("DL44").Value = -2.75938189845475
$ws2.Range("DK45").Value = 0.740740740740751
$ws2.Range("DL45").Value = 3.05676855895196
$ws2.Range("DK46").Value = 3.50421804023362
$ws2.Range("DL46").Value = 5.89459084604716
$ws2.Range("DK47").Value = 2.7863777089783
$ws2.Range("DL47").Value = 1.08303249097472
$ws2.Range("DK48").Value = 1.43160127253447
$ws2.Range("DL48").Value = 4.09429280397024
$ws2.Range("DK49").Value = 1.93887914320009
$ws2.Range("DL49").Value = 0.301175556203236
$ws2.Range("DK50").Value = 1.69491525423729
$ws2.Range("DL50").Value = 1.89431704885345
$ws2.Range("DK51").Value = 8.67924528301886
$ws2.Range("DL51").Value = 10.3448275862069
$ws2.Range("DK52").Value = 0.408163265306138
$ws2.Range("DL52").Value = -0.229182582123746
$ws2.Range("DK53").Value = 3.63712374581941
$ws2.Range("DL53").Value = 4.74922325787838
$ws2.Range("DK54").Value = -1.4519056261343
$ws2.Range("DL54").Value = -0.799999999999997
$ws2.Range("DK55").Value = 1.74767886400875
$ws2.Range("DL55").Value = -0.772430184194897
$ws2.Range("DK56").Value = -2.15053763440859
$ws2.Range("DL56").Value = 0.413223140495874
